# Generate Report for Handback
# - Marks zh-cn / de-de rows as handed back (status text + target/handback
#   file columns + handback datetime) and widens a few columns that now
#   hold longer content.

$wb = $excel.ActiveWorkbook

$srcMdName   = "67ca1668-05ae-4bbd-a9ae-6c043f570a87.md"
$srcMdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57d512d28145d03a5fd9e4934fcd37d128094ff2/e2e/67ca1668-05ae-4bbd-a9ae-6c043f570a87.md"
$statusText  = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status column on every sheet now reads "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsZhCn.Range("C2").Value = $statusText
$wsDeDe.Range("C2").Value = $statusText

# --- zh-cn: record the round-tripped target/handback file + datetime
$wsZhCn.Range("J2").Value = "67ca1668-05ae-4bbd-a9ae-6c043f570a87.bed7c293627cf97d4932911678118c2d963683d3.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-25 22:57:28"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $srcMdUrl, "", "", $srcMdName)

# --- de-de: record the round-tripped target/handback file + datetime
$wsDeDe.Range("J2").Value = "67ca1668-05ae-4bbd-a9ae-6c043f570a87.bed7c293627cf97d4932911678118c2d963683d3.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-25 22:57:35"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $srcMdUrl, "", "", $srcMdName)

# --- Widen columns that now hold the longer status text / hyperlinked file names.
$wsOverview.Range("E1").ColumnWidth = 29.2
$wsOverview.Range("F1").ColumnWidth = 29.2

$wsZhCn.Range("C1").ColumnWidth = 29.2
$wsZhCn.Range("I1").ColumnWidth = 39.15
$wsZhCn.Range("J1").ColumnWidth = 39.15

$wsDeDe.Range("C1").ColumnWidth = 29.2
$wsDeDe.Range("I1").ColumnWidth = 39.15
$wsDeDe.Range("J1").ColumnWidth = 39.15
